# Auto-generated PowerShell Excel COM-interop script
# Applies scheduled market-data refresh values to the Leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 441.17648
$ws.Range("J12").Value = 99.5
$ws.Range("L12").Value = 99.5
$ws.Range("N12").Value = -439.5
$ws.Range("H33").Value = 920.8570999999999
$ws.Range("I33").Value = 699.5454999999999
$ws.Range("K33").Value = 699.5454999999999
$ws.Range("M33").Value = -470.5454999999999
$ws.Range("H70").Value = 5180.815
$ws.Range("J70").Value = 5482.4287
$ws.Range("L70").Value = 16447.2861
$ws.Range("N70").Value = -16987.2861
$ws.Range("H73").Value = 5180.815
$ws.Range("J73").Value = 5482.4287
$ws.Range("L73").Value = 16447.2861
$ws.Range("N73").Value = -18319.2861
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").ClearContents()
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 0
$ws.Range("H137").Value = 12726.844
$ws.Range("I137").Value = 15981.272
$ws.Range("K137").Value = 47943.81600000001
$ws.Range("M137").Value = -45393.81600000001
$ws.Range("H138").Value = 3697.9443
$ws.Range("J138").Value = 4227.5884
$ws.Range("L138").Value = 12682.7652
$ws.Range("N138").Value = -22962.7652

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19126.436
$ws.Range("I32").Value = 21604.982
$ws.Range("J32").Value = 2396.25
$ws.Range("K32").Value = 21604.982
$ws.Range("L32").Value = 2396.25
$ws.Range("M32").Value = -21317.982
$ws.Range("N32").Value = -2970.25
$ws.Range("H61").Value = 4491.1353
$ws.Range("I61").Value = 1450.9166
$ws.Range("J61").Value = 10103.846
$ws.Range("K61").Value = 1450.9166
$ws.Range("L61").Value = 10103.846
$ws.Range("M61").Value = -1238.9166
$ws.Range("N61").Value = -10527.846
$ws.Range("H136").Value = 4491.1353
$ws.Range("I136").Value = 1450.9166
$ws.Range("J136").Value = 10103.846
$ws.Range("K136").Value = 4352.7498
$ws.Range("L136").Value = 30311.538
$ws.Range("M136").Value = -1802.7498
$ws.Range("N136").Value = -35411.538

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1619.25
$ws.Range("I99").Value = 1492.3334
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1492.3334
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = 5.666600000000017
$ws.Range("N99").Value = -4996
$ws.Range("H107").Value = 52928.8
$ws.Range("J107").Value = 3739.8572
$ws.Range("L107").Value = 3739.8572
$ws.Range("N107").Value = -7579.8572
$ws.Range("H134").Value = 2851.423
$ws.Range("I134").Value = 2152.0476
$ws.Range("K134").Value = 6456.1428
$ws.Range("M134").Value = -3921.1428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1671.8572
$ws.Range("J105").Value = 2855.5
$ws.Range("L105").Value = 2855.5
$ws.Range("N105").Value = -6349.5
$ws.Range("H122").Value = 3369.1333
$ws.Range("I122").Value = 3353.25
$ws.Range("J122").Value = 3432.6667
$ws.Range("K122").Value = 10059.75
$ws.Range("L122").Value = 10298.0001
$ws.Range("M122").Value = -7609.75
$ws.Range("N122").Value = -15198.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2206.4443
$ws.Range("I5").Value = 2281.8572
$ws.Range("K5").Value = 6845.571599999999
$ws.Range("M5").Value = -6733.571599999999
$ws.Range("H26").Value = 584.75
$ws.Range("J26").Value = 650
$ws.Range("L26").Value = 1950
$ws.Range("N26").Value = -2526
$ws.Range("H37").Value = 42089.176
$ws.Range("J37").Value = 42089.176
$ws.Range("L37").Value = 126267.528
$ws.Range("N37").Value = -126491.528
$ws.Range("H38").Value = 62500200
$ws.Range("I38").Value = 31.833334
$ws.Range("J38").Value = 100000296
$ws.Range("K38").Value = 95.50000199999999
$ws.Range("L38").Value = 300000888
$ws.Range("M38").Value = 251.499998
$ws.Range("N38").Value = -300001582
$ws.Range("H63").Value = 17958.334
$ws.Range("I63").Value = 17958.334
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 53875.00199999999
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -53126.00199999999
$ws.Range("H66").Value = 17958.334
$ws.Range("I66").Value = 17958.334
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 161625.006
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -157881.006
$ws.Range("H135").Value = 2206.4443
$ws.Range("I135").Value = 2281.8572
$ws.Range("K135").Value = 20536.7148
$ws.Range("M135").Value = -18001.7148
$ws.Range("H141").Value = 4225.1763
$ws.Range("I141").Value = 3796.5454
$ws.Range("J141").Value = 5011
$ws.Range("K141").Value = 11389.6362
$ws.Range("L141").Value = 15033
$ws.Range("M141").Value = -6209.636200000001
$ws.Range("N141").Value = -25393

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 35748.25
$ws.Range("I46").Value = 22995
$ws.Range("J46").Value = 39999.332
$ws.Range("K46").Value = 22995
$ws.Range("L46").Value = 39999.332
$ws.Range("M46").Value = -22839
$ws.Range("N46").Value = -40311.332
$ws.Range("H80").Value = 12011.182
$ws.Range("I80").Value = 6199.25
$ws.Range("J80").Value = 15332.286
$ws.Range("K80").Value = 6199.25
$ws.Range("L80").Value = 15332.286
$ws.Range("M80").Value = -5201.25
$ws.Range("N80").Value = -17328.286
$ws.Range("H83").Value = 12011.182
$ws.Range("I83").Value = 6199.25
$ws.Range("J83").Value = 15332.286
$ws.Range("K83").Value = 30996.25
$ws.Range("L83").Value = 76661.42999999999
$ws.Range("M83").Value = -26004.25
$ws.Range("N83").Value = -86645.42999999999
$ws.Range("H97").Value = 2116.6667
$ws.Range("I97").Value = 1410.7693
$ws.Range("K97").Value = 1410.7693
$ws.Range("M97").Value = -914.7692999999999
$ws.Range("H107").Value = 323
$ws.Range("I107").Value = 258.83334
$ws.Range("J107").Value = 451.33334
$ws.Range("K107").Value = 258.83334
$ws.Range("L107").Value = 451.33334
$ws.Range("M107").Value = 1661.16666
$ws.Range("N107").Value = -4291.33334
$ws.Range("H122").Value = 2702.1
$ws.Range("I122").Value = 2375.5293
$ws.Range("K122").Value = 7126.5879
$ws.Range("M122").Value = -4676.5879
$ws.Range("H126").Value = 3480.087
$ws.Range("I126").Value = 1990.091
$ws.Range("J126").Value = 4845.9165
$ws.Range("K126").Value = 5970.272999999999
$ws.Range("L126").Value = 14537.7495
$ws.Range("M126").Value = -3500.272999999999
$ws.Range("N126").Value = -19477.7495
$ws.Range("H132").Value = 2493.8936
$ws.Range("I132").Value = 2508.9783
$ws.Range("K132").Value = 7526.9349
$ws.Range("M132").Value = -4996.9349

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3580.3076
$ws.Range("I46").Value = 1091.875
$ws.Range("J46").Value = 7561.8
$ws.Range("K46").Value = 1091.875
$ws.Range("L46").Value = 7561.8
$ws.Range("M46").Value = -903.875
$ws.Range("N46").Value = -7937.8
$ws.Range("H48").Value = 120000
$ws.Range("J48").Value = 120000
$ws.Range("L48").Value = 120000
$ws.Range("N48").Value = -121322
$ws.Range("H55").Value = 882.8261
$ws.Range("J55").Value = 1275.3077
$ws.Range("L55").Value = 1275.3077
$ws.Range("N55").Value = -1621.3077
$ws.Range("H132").Value = 2483.4062
$ws.Range("I132").Value = 1859.16
$ws.Range("K132").Value = 5577.48
$ws.Range("M132").Value = -3047.48

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 23351.25
$ws.Range("I38").Value = 2000
$ws.Range("J38").Value = 30468.334
$ws.Range("K38").Value = 2000
$ws.Range("L38").Value = 30468.334
$ws.Range("M38").Value = -1527
$ws.Range("N38").Value = -31414.334
$ws.Range("H63").Value = 47499.5
$ws.Range("J63").Value = 47499.5
$ws.Range("L63").Value = 47499.5
$ws.Range("N63").Value = -48747.5
$ws.Range("H66").Value = 47499.5
$ws.Range("J66").Value = 47499.5
$ws.Range("L66").Value = 142498.5
$ws.Range("N66").Value = -148738.5
